$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the newest Adafruit IO reading as a new row, matching the shape
# and style of the existing data rows.
$newRow = @("2024-09-25T18:06:40Z", "temperature", "25", "N/A", "N/A", "N/A")
$rowIndex = 57

for ($col = 1; $col -le 6; $col++) {
    $cell = $ws.Cells.Item($rowIndex, $col)
    # Prefix with an apostrophe so numeric-looking text (e.g. "25") is
    # stored as text, consistent with the rest of the sheet.
    $cell.Value = "'" + $newRow[$col - 1]
}

# Copy the formatting from the row above so the new row keeps the same
# (default) cell style instead of picking up a text/quote-prefix style.
$src = $ws.Range("A56:F56")
$dst = $ws.Range("A57:F57")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
